$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 04:55:26"
$wsZhCn.Range("H2").Value = "2016-03-19 04:55:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 04:55:29"
$wsDeDe.Range("H2").Value = "2016-03-19 04:55:45"
